# Apply the edits described by the diff for Data/Per RHA/Temperature.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the temperature values in row 9 (VIII - Lerma Santiago Pacifico) ---
$ws.Range("C9").Value = 19.581270816127027
$ws.Range("D9").Value = 19.422435795585798
$ws.Range("E9").Value = 20.156168623234638
$ws.Range("F9").Value = 20.598878795780386
$ws.Range("G9").Value = 20.292354995371639
$ws.Range("H9").Value = 20.077527466934161
$ws.Range("I9").Value = 20.2794381352588
$ws.Range("J9").Value = 19.415551616304157
$ws.Range("K9").Value = 20.023270844569232
$ws.Range("L9").Value = 19.969374999999999
$ws.Range("M9").Value = 20.400215618566513
$ws.Range("N9").Value = 20.338541666666668
$ws.Range("O9").Value = 20.599999999999998
$ws.Range("P9").Value = 20.689583333333335
$ws.Range("Q9").Value = 20.681249999999999
$ws.Range("R9").Value = 20.664583333333333
$ws.Range("S9").Value = 21.174999999999997
$ws.Range("T9").Value = 20.837500000000002

# --- Update the view state: scroll the window so column I is at the left edge
#     and select cell Q11 (matching the sheetView's topLeftCell/selection) ---
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 9
$aw.ScrollRow = 1
$null = $ws.Range("Q11").Select()
